# Add seven new monster rows (532-538) for the two new passive-skill quests
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row9 = New-Object "object[,]" 1,50
$row9[0,0] = 532
$row9[0,1] = "Duke of Suicide"
$row9[0,2] = 100000000000
$row9[0,3] = 100000000000
$row9[0,4] = 100000000000
$row9[0,5] = 100000000000
$row9[0,6] = 100000000000
$row9[0,7] = 100000000000
$row9[0,8] = 100000000000
$row9[0,9] = 50000000000
$row9[0,10] = 1.1
$row9[0,11] = 1.1
$row9[0,12] = 1.1
$row9[0,13] = 1.1
$row9[0,15] = 0
$row9[0,16] = 0
$row9[0,17] = 1
$row9[0,18] = 999
$row9[0,19] = "str"
$row9[0,20] = 200000
$row9[0,21] = 0.1
$row9[0,22] = 2000000000000
$row9[0,23] = 0
$row9[0,24] = "200000000000-250000000000"
$row9[0,25] = "150000000000-225000000000"
$row9[0,26] = 25000000000
$row9[0,27] = 10000000000
$row9[0,28] = 1.5
$row9[0,29] = 1.25
$row9[0,30] = 1.25
$row9[0,31] = 1.25
$row9[0,32] = 0
$row9[0,33] = 0
$row9[0,34] = 1.05
$row9[0,35] = 1.05
$row9[0,36] = 1.05
$row9[0,37] = 1.05
$row9[0,39] = 0
$row9[0,40] = "Delusional Memories"
$row9[0,42] = 0.7
$row9[0,43] = 0.45
$row9[0,44] = 0.35
$row9[0,48] = 0.15
$row9[0,49] = 9
$ws.Range("A9:AX9").Value = $row9

$row10 = New-Object "object[,]" 1,50
$row10[0,0] = 533
$row10[0,1] = "Lumbering Lord of Nightmares"
$row10[0,2] = 100000000000
$row10[0,3] = 100000000000
$row10[0,4] = 100000000000
$row10[0,5] = 100000000000
$row10[0,6] = 100000000000
$row10[0,7] = 100000000000
$row10[0,8] = 100000000000
$row10[0,9] = 50000000000
$row10[0,10] = 1.1
$row10[0,11] = 1.1
$row10[0,12] = 1.1
$row10[0,13] = 1.1
$row10[0,15] = 0
$row10[0,16] = 0
$row10[0,17] = 1
$row10[0,18] = 999
$row10[0,19] = "str"
$row10[0,20] = 200000
$row10[0,21] = 0.1
$row10[0,22] = 2000000000000
$row10[0,23] = 0
$row10[0,24] = "200000000000-250000000000"
$row10[0,25] = "150000000000-225000000000"
$row10[0,26] = 25000000000
$row10[0,27] = 10000000000
$row10[0,28] = 1.5
$row10[0,29] = 1.25
$row10[0,30] = 1.25
$row10[0,31] = 1.25
$row10[0,32] = 0
$row10[0,33] = 0
$row10[0,34] = 1.05
$row10[0,35] = 1.05
$row10[0,36] = 1.05
$row10[0,37] = 1.05
$row10[0,39] = 0
$row10[0,40] = "Delusional Memories"
$row10[0,42] = 0.7
$row10[0,43] = 0.45
$row10[0,44] = 0.35
$row10[0,48] = 0.15
$row10[0,49] = 9
$ws.Range("A10:AX10").Value = $row10

$row11 = New-Object "object[,]" 1,50
$row11[0,0] = 534
$row11[0,1] = "Shadow Lord of Children's Tears"
$row11[0,2] = 100000000000
$row11[0,3] = 100000000000
$row11[0,4] = 100000000000
$row11[0,5] = 100000000000
$row11[0,6] = 100000000000
$row11[0,7] = 100000000000
$row11[0,8] = 100000000000
$row11[0,9] = 50000000000
$row11[0,10] = 1.1
$row11[0,11] = 1.1
$row11[0,12] = 1.1
$row11[0,13] = 1.1
$row11[0,15] = 0
$row11[0,16] = 0
$row11[0,17] = 1
$row11[0,18] = 999
$row11[0,19] = "str"
$row11[0,20] = 200000
$row11[0,21] = 0.1
$row11[0,22] = 2000000000000
$row11[0,23] = 0
$row11[0,24] = "200000000000-250000000000"
$row11[0,25] = "150000000000-225000000000"
$row11[0,26] = 25000000000
$row11[0,27] = 10000000000
$row11[0,28] = 1.5
$row11[0,29] = 1.25
$row11[0,30] = 1.25
$row11[0,31] = 1.25
$row11[0,32] = 0
$row11[0,33] = 0
$row11[0,34] = 1.05
$row11[0,35] = 1.05
$row11[0,36] = 1.05
$row11[0,37] = 1.05
$row11[0,39] = 0.15
$row11[0,40] = "Delusional Memories"
$row11[0,42] = 0.7
$row11[0,43] = 0.45
$row11[0,44] = 0.35
$row11[0,48] = 0.15
$row11[0,49] = 9
$ws.Range("A11:AX11").Value = $row11

$row12 = New-Object "object[,]" 1,50
$row12[0,0] = 535
$row12[0,1] = "Princess of Dead Cats Wishes"
$row12[0,2] = 100000000000
$row12[0,3] = 100000000000
$row12[0,4] = 100000000000
$row12[0,5] = 100000000000
$row12[0,6] = 100000000000
$row12[0,7] = 100000000000
$row12[0,8] = 100000000000
$row12[0,9] = 50000000000
$row12[0,10] = 1.1
$row12[0,11] = 1.1
$row12[0,12] = 1.1
$row12[0,13] = 1.1
$row12[0,15] = 0
$row12[0,16] = 0
$row12[0,17] = 1
$row12[0,18] = 999
$row12[0,19] = "str"
$row12[0,20] = 200000
$row12[0,21] = 0.1
$row12[0,22] = 2000000000000
$row12[0,23] = 0
$row12[0,24] = "200000000000-250000000000"
$row12[0,25] = "150000000000-225000000000"
$row12[0,26] = 25000000000
$row12[0,27] = 10000000000
$row12[0,28] = 1.5
$row12[0,29] = 1.25
$row12[0,30] = 1.25
$row12[0,31] = 1.25
$row12[0,32] = 0
$row12[0,33] = 0
$row12[0,34] = 1.05
$row12[0,35] = 1.05
$row12[0,36] = 1.05
$row12[0,37] = 1.05
$row12[0,39] = 0.15
$row12[0,40] = "Delusional Memories"
$row12[0,42] = 0.7
$row12[0,43] = 0.45
$row12[0,44] = 0.35
$row12[0,48] = 0.15
$row12[0,49] = 9
$ws.Range("A12:AX12").Value = $row12

$row13 = New-Object "object[,]" 1,50
$row13[0,0] = 536
$row13[0,1] = "Time Jester of Satan"
$row13[0,2] = 100000000000
$row13[0,3] = 100000000000
$row13[0,4] = 100000000000
$row13[0,5] = 100000000000
$row13[0,6] = 100000000000
$row13[0,7] = 100000000000
$row13[0,8] = 100000000000
$row13[0,9] = 50000000000
$row13[0,10] = 1.1
$row13[0,11] = 1.1
$row13[0,12] = 1.1
$row13[0,13] = 1.1
$row13[0,15] = 0
$row13[0,16] = 0
$row13[0,17] = 1
$row13[0,18] = 999
$row13[0,19] = "str"
$row13[0,20] = 200000
$row13[0,21] = 0.1
$row13[0,22] = 2000000000000
$row13[0,23] = 0
$row13[0,24] = "200000000000-250000000000"
$row13[0,25] = "150000000000-225000000000"
$row13[0,26] = 25000000000
$row13[0,27] = 10000000000
$row13[0,28] = 1.5
$row13[0,29] = 1.25
$row13[0,30] = 1.25
$row13[0,31] = 1.25
$row13[0,32] = 0
$row13[0,33] = 0
$row13[0,34] = 1.05
$row13[0,35] = 1.05
$row13[0,36] = 1.05
$row13[0,37] = 1.05
$row13[0,39] = 0.15
$row13[0,40] = "Delusional Memories"
$row13[0,42] = 0.7
$row13[0,43] = 0.45
$row13[0,44] = 0.35
$row13[0,48] = 0.15
$row13[0,49] = 9
$ws.Range("A13:AX13").Value = $row13

$row14 = New-Object "object[,]" 1,50
$row14[0,0] = 537
$row14[0,1] = "Satan"
$row14[0,2] = 100000000000
$row14[0,3] = 100000000000
$row14[0,4] = 100000000000
$row14[0,5] = 100000000000
$row14[0,6] = 100000000000
$row14[0,7] = 100000000000
$row14[0,8] = 100000000000
$row14[0,9] = 50000000000
$row14[0,10] = 1.1
$row14[0,11] = 1.1
$row14[0,12] = 1.1
$row14[0,13] = 1.1
$row14[0,15] = 0
$row14[0,16] = 0
$row14[0,17] = 1
$row14[0,18] = 999
$row14[0,19] = "str"
$row14[0,20] = 200000
$row14[0,21] = 0.1
$row14[0,22] = 2000000000000
$row14[0,23] = 0
$row14[0,24] = "200000000000-250000000000"
$row14[0,25] = "150000000000-225000000000"
$row14[0,26] = 25000000000
$row14[0,27] = 10000000000
$row14[0,28] = 1.5
$row14[0,29] = 1.25
$row14[0,30] = 1.25
$row14[0,31] = 1.25
$row14[0,32] = 0
$row14[0,33] = 0
$row14[0,34] = 1.05
$row14[0,35] = 1.05
$row14[0,36] = 1.05
$row14[0,37] = 1.05
$row14[0,39] = 0.15
$row14[0,40] = "Delusional Memories"
$row14[0,42] = 0.7
$row14[0,43] = 0.45
$row14[0,44] = 0.35
$row14[0,48] = 0.15
$row14[0,49] = 9
$ws.Range("A14:AX14").Value = $row14

$row15 = New-Object "object[,]" 1,50
$row15[0,0] = 538
$row15[0,1] = "Corrupted Satan of The Deep"
$row15[0,2] = 100000000000
$row15[0,3] = 100000000000
$row15[0,4] = 100000000000
$row15[0,5] = 100000000000
$row15[0,6] = 100000000000
$row15[0,7] = 100000000000
$row15[0,8] = 100000000000
$row15[0,9] = 50000000000
$row15[0,10] = 1.1
$row15[0,11] = 1.1
$row15[0,12] = 1.1
$row15[0,13] = 1.1
$row15[0,15] = 0
$row15[0,16] = 0
$row15[0,17] = 1
$row15[0,18] = 999
$row15[0,19] = "str"
$row15[0,20] = 200000
$row15[0,21] = 0.1
$row15[0,22] = 2000000000000
$row15[0,23] = 0
$row15[0,24] = "200000000000-250000000000"
$row15[0,25] = "150000000000-225000000000"
$row15[0,26] = 25000000000
$row15[0,27] = 10000000000
$row15[0,28] = 1.5
$row15[0,29] = 1.25
$row15[0,30] = 1.25
$row15[0,31] = 1.25
$row15[0,32] = 0
$row15[0,33] = 0
$row15[0,34] = 1.05
$row15[0,35] = 1.05
$row15[0,36] = 1.05
$row15[0,37] = 1.05
$row15[0,38] = "Broken Blessed Sword of Pags"
$row15[0,39] = 0.14
$row15[0,40] = "Delusional Memories"
$row15[0,42] = 0.7
$row15[0,43] = 0.45
$row15[0,44] = 0.35
$row15[0,48] = 0.15
$row15[0,49] = 9
$ws.Range("A15:AX15").Value = $row15

# Widen columns B (name) and AM (quest_item_name) to fit the new, longer text
$ws.Columns.Item(2).ColumnWidth = 36.83
$ws.Columns.Item(39).ColumnWidth = 33.33

